$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "72.032.88"
$ws.Range("E2").Value = "  +3.51%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "4.046.76"
$ws.Range("E3").Value = "  +2.95%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.00"

# Row 5 - BNB
$ws.Range("D5").Value = "'524.61"
$ws.Range("E5").Value = "  -2.12%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'148.98"
$ws.Range("E6").Value = "  +2.59%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = "  +1.26%  "

# Row 8 - USDC
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.12%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.71%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.177"
$ws.Range("E10").Value = "  +2.28%  "

# Row 11 - ShibaInu
$ws.Range("D11").Value = "'0.0000342"
$ws.Range("E11").Value = "  +0.82%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'46.86"
$ws.Range("E12").Value = "  +9.93%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'10.80"
$ws.Range("E13").Value = "  +4.02%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.683.64"
$ws.Range("E14").Value = "  +2.67%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "4.039.10"
$ws.Range("E15").Value = "  +2.63%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "'21.54"
$ws.Range("E16").Value = "  +8.67%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "'14.37"
$ws.Range("E17").Value = "  +2.54%  "

# Row 18 - Polygon
$ws.Range("E18").Value = "  +0.49%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  -1.71%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "71.981.49"
$ws.Range("E20").Value = "  +3.62%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'443.96"
$ws.Range("E21").Value = "  +2.73%  "

# Row 22 - ImmutableX
$ws.Range("E22").Value = "  +5.70%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "'94.83"
$ws.Range("E23").Value = "  +6.95%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "'14.38"
$ws.Range("E24").Value = "  -0.92%  "

# Row 25 - RenderToken
$ws.Range("D25").Value = "'12.37"
$ws.Range("E25").Value = "  +4.83%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "'4.05"
$ws.Range("E26").Value = "  -3.86%  "

# Row 27 - Filecoin
$ws.Range("E27").Value = "  +3.22%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'37.33"
$ws.Range("E28").Value = "  +1.85%  "

# Row 29 - Cosmos
$ws.Range("E29").Value = "  +2.28%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "'700.78"
$ws.Range("E30").Value = "  -0.01%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +3.24%  "

# Row 32 - Toncoin
$ws.Range("E32").Value = "  +2.43%  "

# Row 33 - NEARProtocol
$ws.Range("E33").Value = "  +13.59%  "

# Row 34 - OKB
$ws.Range("D34").Value = "'68.25"
$ws.Range("E34").Value = "  -6.52%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "0.0₃0917"
$ws.Range("E35").Value = "  +6.28%  "

# Row 36 - TheGraph
$ws.Range("E36").Value = "  -3.80%  "

# Row 37 - InjectiveProtocol
$ws.Range("D37").Value = "'40.94"
$ws.Range("E37").Value = "  +0.96%  "

# Row 38 - Kaspa
$ws.Range("E38").Value = "  +5.66%  "

# Row 39 - ThetaToken
$ws.Range("D39").Value = "'3.55"
$ws.Range("E39").Value = "  +18.46%  "

# Row 40 - Dai
$ws.Range("E40").Value = "  +0.26%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  -0.23%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +1.36%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  +0.66%  "

# Row 44 - Fetch.AI
$ws.Range("E44").Value = "  -0.62%  "

# Row 45 - ApeXProtocol
$ws.Range("D45").Value = "'3.54"
$ws.Range("E45").Value = "  +4.57%  "

# Row 46 - was Stellar -> now Stacks (rows 46/47 swapped)
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'3.23"
$ws.Range("E46").Value = "  +1.49%  "

# Row 47 - was Stacks -> now Stellar (rows 46/47 swapped)
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").Value = "'0.146"
$ws.Range("E47").Value = "  +2.45%  "

# Row 48 - FLOKI
$ws.Range("D48").Value = "'0.000281"
$ws.Range("E48").Value = "  +18.94%  "

# Row 49 - THORChain
$ws.Range("D49").Value = "'9.22"
$ws.Range("E49").Value = "  +6.07%  "

# Row 50 - LidoDAOToken
$ws.Range("E50").Value = "  +1.43%  "

# Row 51 - BabyDogeCoin
$ws.Range("E51").Value = "  -1.51%  "
